$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.289.72"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "'2.283.58"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'318.60"
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("D6").Value = "'100.93"
$ws.Range("E6").Value = "  -4.67%  "
$ws.Range("D7").Value = "'0.626"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "'0.601"
$ws.Range("E9").Value = "  -1.44%  "
$ws.Range("D10").Value = "'39.06"
$ws.Range("E10").Value = "  -2.31%  "
$ws.Range("D11").Value = "'0.0901"
$ws.Range("E11").Value = "  -1.07%  "
$ws.Range("D12").Value = "'8.22"
$ws.Range("E12").Value = "  -1.90%  "
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").Value = "'0.956"
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").Value = "'15.13"
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").Value = "'2.631.11"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "'2.288.71"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "'42.263.98"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").Value = "'7.37"
$ws.Range("E19").Value = "  -3.36%  "
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D21").Value = "'12.97"
$ws.Range("E21").Value = "  +32.80%  "
$ws.Range("D22").Value = "'72.67"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("D24").Value = "'266.50"
$ws.Range("E24").Value = "  +2.75%  "
$ws.Range("D25").Value = "'2.20"
$ws.Range("E25").Value = "  -5.28%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'10.78"
$ws.Range("E28").Value = "  +2.58%  "
$ws.Range("D29").Value = "'22.42"
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("D30").Value = "'37.26"
$ws.Range("E30").Value = "  +2.67%  "
$ws.Range("D31").Value = "'165.76"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "'6.09"
$ws.Range("E32").Value = "  +3.30%  "
$ws.Range("D33").Value = "'0.0869"
$ws.Range("E33").Value = "  -2.38%  "
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("E35").Value = "  -12.72%  "
$ws.Range("E36").Value = "  -3.42%  "
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("D38").Value = "'0.0356"
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "'2.74"
$ws.Range("E39").Value = "  -6.15%  "
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").Value = "'3.63"
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("E41").Value = "  +2.35%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "'68.39"
$ws.Range("E43").Value = "  -3.56%  "
$ws.Range("D44").Value = "'0.223"
$ws.Range("E44").Value = "  -1.80%  "
$ws.Range("D45").Value = "'92.84"
$ws.Range("E45").Value = "  -6.69%  "
$ws.Range("D46").Value = "'114.44"
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("D47").Value = "'11.88"
$ws.Range("E47").Value = "  -2.61%  "
$ws.Range("D48").Value = "'78.60"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").Value = "'8.93"
$ws.Range("E49").Value = "  -2.29%  "
$ws.Range("D50").Value = "'5.20"
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("D51").Value = "'1.597.23"
$ws.Range("E51").Value = "  +3.28%  "
